# SSDM-55: fixed more tests related to property creation and xls export.
#
# Adds a new "Multivalued" column (K) to the sample-type export sheet:
#   - K4 gets the new column header "Multivalued" (bold, 14pt, black, Calibri)
#   - K5:K9 get the text value "FALSE" (same boolean-like text style used
#     elsewhere in the sheet, e.g. column B)
#   - Row 4's height shrinks slightly (18.75 -> 17.35)
#   - Selection moves to K12 (last used cell after the edit)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- New column header: K4 = "Multivalued" ------------------------------
$header = $ws.Range("K4")
$header.Value = "Multivalued"
$header.Font.Name = "Calibri"
$header.Font.Bold = $true
$header.Font.Size = 14
$header.Font.Color = 0

# --- New column values: K5:K9 = "FALSE" ----------------------------------
# Leading apostrophe forces these to be stored as text (matching the
# existing text-typed "FALSE"/"TRUE" cells elsewhere in the sheet) instead
# of being auto-converted to native booleans.
$boolLikeFormat = $ws.Range("B5").NumberFormat
foreach ($r in 5..9) {
    $cell = $ws.Range("K" + $r)
    $cell.Value = "'FALSE"
    $cell.NumberFormat = $boolLikeFormat
    $cell.HorizontalAlignment = $ws.Range("B5").HorizontalAlignment
}

# --- Row 4 height adjustment ---------------------------------------------
$ws.Rows.Item(4).RowHeight = 17.35

# --- Final selection -------------------------------------------------------
$null = $ws.Range("K12").Select()
